# Updates cryptos list values (price/volume) and reorders a few coin rows
# to match the latest data pull, per commit message:
# "Updated cryptos list on Tue Nov 21 20:52:26 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.069.99'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.12%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.996.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.77%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.606'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.87%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("E8").Value = '  -3.78%  '

$ws.Range("E9").Value = '  -3.45%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.07'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.37%  '

$ws.Range("E11").Value = '  -4.92%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0981'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.02%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.295.54'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.58%  '

$ws.Range("E14").Value = '  -4.30%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.95'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.32%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.762'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.013.97'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '37.046.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.32%  '

$ws.Range("E21").Value = '  -4.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.46%  '

$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.35'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.125'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.83%  '

$ws.Range("E31").Value = '  -3.21%  '

$ws.Range("E32").Value = '  -2.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.44'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.21%  '

$ws.Range("E34").Value = '  -8.29%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.25'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.53%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.20%  '

$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.85%  '

$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.79'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.63%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.54%  '

$ws.Range("E41").Value = '  +2.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.439.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.52%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.13'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.53%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0204'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.29%  '

$ws.Range("E45").Value = '  -8.74%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.68'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.89%  '

$ws.Range("E48").Value = '  -3.79%  '

$ws.Range("E49").Value = '  +0.69%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.71'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.79%  '

$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.186.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.66%  '
